$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '61.409.40'
$ws.Range("E2").Value = '  -4.64%  '
$ws.Range("D3").Value = '3.302.85'
$ws.Range("E3").Value = '  -5.64%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.10%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '565.60'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -4.22%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '125.87'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -6.41%  '
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("D8").Value = '3.298.57'
$ws.Range("E8").Value = '  -5.72%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.473'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -3.05%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.18'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -5.72%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.117'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -6.36%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.373'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -4.96%  '
$ws.Range("D13").Value = '3.849.94'
$ws.Range("E13").Value = '  -6.02%  '
$ws.Range("E14").Value = '  -1.67%  '
$ws.Range("D15").Value = '3.287.99'
$ws.Range("E15").Value = '  -6.09%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0000168'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -6.99%  '
$ws.Range("D17").Value = '61.462.22'
$ws.Range("E17").Value = '  -4.55%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '24.50'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -4.84%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '9.11'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -9.48%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '5.56'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -4.38%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.07'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -3.82%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '355.21'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -9.25%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.548'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -6.00%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.999'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.18%  '
$ws.Range("D25").Value = '3.422.54'
$ws.Range("E25").Value = '  -5.98%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '69.93'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -6.05%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0000106'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -8.88%  '
$ws.Range("E28").Value = '  +0.01%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.10'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -4.67%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.45'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -2.51%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.84'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -4.50%  '
$ws.Range("B32").Value = 'PancakeSwap'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.09'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -8.16%  '
$ws.Range("B33").Value = 'USDe'
$ws.Range("C33").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.00'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.09%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.147'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -6.08%  '
$ws.Range("D35").Value = '3.319.32'
$ws.Range("E35").Value = '  -5.92%  '
$ws.Range("B36").Value = 'NEARProtocol'
$ws.Range("C36").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.46'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +2.13%  '
$ws.Range("B37").Value = 'EthereumClassic'
$ws.Range("C37").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '22.30'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -4.92%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '166.10'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.06%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '6.62'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -5.19%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.49'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -4.92%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0754'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -4.96%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.998'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.25%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '40.92'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.77%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.745'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -8.04%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '4.16'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -6.51%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.10'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -7.25%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.53'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -8.10%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '22.30'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -11.02%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '6.55'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -3.99%  '
$ws.Range("B50").Value = 'SuiNetwork'
$ws.Range("C50").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.847'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -8.80%  '
$ws.Range("B51").Value = 'Maker'
$ws.Range("C51").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D51").Value = '2.193.92'
$ws.Range("E51").Value = '  -8.97%  '
